$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.273.00"
$ws.Range("E2").Value = "  +4.49%  "
$ws.Range("D3").Value = "3.627.03"
$ws.Range("E3").Value = "  +5.11%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'202.18"
$ws.Range("E5").Value = "  +12.59%  "
$ws.Range("D6").Value = "'578.46"
$ws.Range("E6").Value = "  +4.25%  "
$ws.Range("E7").Value = "  +4.53%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.685"
$ws.Range("E9").Value = "  +7.32%  "
$ws.Range("D10").Value = "'61.37"
$ws.Range("E10").Value = "  +20.87%  "
$ws.Range("D11").Value = "'0.149"
$ws.Range("E11").Value = "  +7.27%  "
$ws.Range("E12").Value = "  +15.01%  "
$ws.Range("D13").Value = "'10.26"
$ws.Range("E13").Value = "  +9.78%  "
$ws.Range("D14").Value = "4.197.30"
$ws.Range("E14").Value = "  +4.95%  "
$ws.Range("D15").Value = "3.624.85"
$ws.Range("D16").Value = "'19.45"
$ws.Range("E16").Value = "  +10.94%  "
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").Value = "68.161.19"
$ws.Range("E18").Value = "  +4.98%  "
$ws.Range("D19").Value = "'12.41"
$ws.Range("E19").Value = "  +7.71%  "
$ws.Range("E20").Value = "  +5.76%  "
$ws.Range("D21").Value = "'408.87"
$ws.Range("E21").Value = "  +9.61%  "
$ws.Range("D22").Value = "'13.15"
$ws.Range("E22").Value = "  +23.77%  "
$ws.Range("D23").Value = "'4.24"
$ws.Range("E23").Value = "  +5.12%  "
$ws.Range("D24").Value = "'85.70"
$ws.Range("E24").Value = "  +4.67%  "
$ws.Range("D25").Value = "'4.00"
$ws.Range("E25").Value = "  +18.57%  "
$ws.Range("D26").Value = "'2.94"
$ws.Range("E26").Value = "  +6.67%  "
$ws.Range("D27").Value = "'12.63"
$ws.Range("E27").Value = "  +8.16%  "
$ws.Range("E28").Value = "  +2.42%  "
$ws.Range("D29").Value = "'9.36"
$ws.Range("E29").Value = "  +11.04%  "
$ws.Range("D30").Value = "'7.83"
$ws.Range("E30").Value = "  +11.30%  "
$ws.Range("D31").Value = "'31.82"
$ws.Range("E31").Value = "  +6.66%  "
$ws.Range("D32").Value = "'683.14"
$ws.Range("E32").Value = "  +12.27%  "
$ws.Range("D33").Value = "'12.26"
$ws.Range("E33").Value = "  +4.79%  "
$ws.Range("E34").Value = "  +6.22%  "
$ws.Range("D35").Value = "'63.92"
$ws.Range("E35").Value = "  +2.48%  "
$ws.Range("D36").Value = "'41.98"
$ws.Range("E36").Value = "  +4.94%  "
$ws.Range("E37").Value = "  +6.38%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").Value = "0.0₃0769"
$ws.Range("E39").Value = "  +9.18%  "
$ws.Range("D40").Value = "'3.20"
$ws.Range("E40").Value = "  +19.69%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "'2.72"
$ws.Range("E44").Value = "  +12.43%  "
$ws.Range("D45").Value = "'2.89"
$ws.Range("E45").Value = "  +29.77%  "
$ws.Range("D46").Value = "'2.87"
$ws.Range("E46").Value = "  +18.13%  "
$ws.Range("E47").Value = "  +7.28%  "
$ws.Range("E48").Value = "  +5.70%  "
$ws.Range("D49").Value = "'8.79"
$ws.Range("E49").Value = "  +9.47%  "
$ws.Range("D50").Value = "'3.08"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").Value = "'139.54"
$ws.Range("E51").Value = "  +1.48%  "

# Row 41/42 content swap (Maker <-> Kaspa) with updated values
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.136"
$ws.Range("E41").Value = "  +6.22%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.194.52"
$ws.Range("E42").Value = "  +10.58%  "
